$d = $word.ActiveDocument
$d.Content.Find.Execute("Social Service Provider", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Internship", 2)
